$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as the first data row
# (row 7), pushing the existing rows 7-16 down to rows 8-17.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record's data.
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44789
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100112017
$ws.Cells.Item(7, 7).Value = "Ramas de apio"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 5000
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = 5500
$ws.Cells.Item(7, 14).Value = "`$/atado 7 kilos"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 5500
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
